$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix material waste index label in reports:
# The header above the summary table (merged cell C8:F8) was labeled
# "Waste Index" but should read "Material Consumption Index".
$ws.Range("C8").Value = "Material Consumption Index"

# Move the active selection from I11 to A10 (matches the author's saved
# selection state in the workbook).
$ws.Range("A10").Select()
